$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.288.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.492.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.39%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.36%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("E8").Value = "  -2.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.491.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.34%  "

$ws.Range("E10").Value = "  -3.56%  "

$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("E12").Value = "  -2.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.970.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.38%  "

$ws.Range("E16").Value = "  -2.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.181.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.486.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.14%  "

$ws.Range("E23").Value = "  -4.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("E26").Value = "  -6.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.89%  "

$ws.Range("E29").Value = "  -0.73%  "

$ws.Range("E30").Value = "  -5.99%  "

$ws.Range("E31").Value = "  -1.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "517.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.65%  "

$ws.Range("E33").Value = "  -1.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.92%  "

$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("E36").Value = "  -2.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.15%  "

$ws.Range("E38").Value = "  +1.44%  "

$ws.Range("E39").Value = "  -2.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.45%  "

$ws.Range("E41").Value = "  -2.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.332"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.90%  "

$ws.Range("E44").Value = "  +0.20%  "

$ws.Range("E45").Value = "  -1.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.15%  "

$ws.Range("E48").Value = "  -2.72%  "

$ws.Range("E49").Value = "  -4.30%  "

$ws.Range("E50").Value = "  -3.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.35%  "
